$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot original values for the columns that get permuted across rows
# (D=Fecha, J=Volumen, K=Precio minimo, L=Precio maximo, M=Precio promedio ponderado, P=Precio $/Kg)
# Use .Value2 (plain numeric) instead of .Value, since .Value getter is unreliable in this runtime.
$orig = @{}
$orig[2] = @{
    D = $ws.Cells.Item(2, 4).Value2
    J = $ws.Cells.Item(2, 10).Value2
    K = $ws.Cells.Item(2, 11).Value2
    L = $ws.Cells.Item(2, 12).Value2
    M = $ws.Cells.Item(2, 13).Value2
    P = $ws.Cells.Item(2, 16).Value2
}
$orig[3] = @{
    D = $ws.Cells.Item(3, 4).Value2
    J = $ws.Cells.Item(3, 10).Value2
    K = $ws.Cells.Item(3, 11).Value2
    L = $ws.Cells.Item(3, 12).Value2
    M = $ws.Cells.Item(3, 13).Value2
    P = $ws.Cells.Item(3, 16).Value2
}
$orig[4] = @{
    D = $ws.Cells.Item(4, 4).Value2
    J = $ws.Cells.Item(4, 10).Value2
    K = $ws.Cells.Item(4, 11).Value2
    L = $ws.Cells.Item(4, 12).Value2
    M = $ws.Cells.Item(4, 13).Value2
    P = $ws.Cells.Item(4, 16).Value2
}
$orig[5] = @{
    D = $ws.Cells.Item(5, 4).Value2
    J = $ws.Cells.Item(5, 10).Value2
    K = $ws.Cells.Item(5, 11).Value2
    L = $ws.Cells.Item(5, 12).Value2
    M = $ws.Cells.Item(5, 13).Value2
    P = $ws.Cells.Item(5, 16).Value2
}
$orig[6] = @{
    D = $ws.Cells.Item(6, 4).Value2
    J = $ws.Cells.Item(6, 10).Value2
    K = $ws.Cells.Item(6, 11).Value2
    L = $ws.Cells.Item(6, 12).Value2
    M = $ws.Cells.Item(6, 13).Value2
    P = $ws.Cells.Item(6, 16).Value2
}
$orig[7] = @{
    D = $ws.Cells.Item(7, 4).Value2
    J = $ws.Cells.Item(7, 10).Value2
    K = $ws.Cells.Item(7, 11).Value2
    L = $ws.Cells.Item(7, 12).Value2
    M = $ws.Cells.Item(7, 13).Value2
    P = $ws.Cells.Item(7, 16).Value2
}
$orig[8] = @{
    D = $ws.Cells.Item(8, 4).Value2
    J = $ws.Cells.Item(8, 10).Value2
    K = $ws.Cells.Item(8, 11).Value2
    L = $ws.Cells.Item(8, 12).Value2
    M = $ws.Cells.Item(8, 13).Value2
    P = $ws.Cells.Item(8, 16).Value2
}
$orig[9] = @{
    D = $ws.Cells.Item(9, 4).Value2
    J = $ws.Cells.Item(9, 10).Value2
    K = $ws.Cells.Item(9, 11).Value2
    L = $ws.Cells.Item(9, 12).Value2
    M = $ws.Cells.Item(9, 13).Value2
    P = $ws.Cells.Item(9, 16).Value2
}
$orig[10] = @{
    D = $ws.Cells.Item(10, 4).Value2
    J = $ws.Cells.Item(10, 10).Value2
    K = $ws.Cells.Item(10, 11).Value2
    L = $ws.Cells.Item(10, 12).Value2
    M = $ws.Cells.Item(10, 13).Value2
    P = $ws.Cells.Item(10, 16).Value2
}
$orig[11] = @{
    D = $ws.Cells.Item(11, 4).Value2
    J = $ws.Cells.Item(11, 10).Value2
    K = $ws.Cells.Item(11, 11).Value2
    L = $ws.Cells.Item(11, 12).Value2
    M = $ws.Cells.Item(11, 13).Value2
    P = $ws.Cells.Item(11, 16).Value2
}
$orig[12] = @{
    D = $ws.Cells.Item(12, 4).Value2
    J = $ws.Cells.Item(12, 10).Value2
    K = $ws.Cells.Item(12, 11).Value2
    L = $ws.Cells.Item(12, 12).Value2
    M = $ws.Cells.Item(12, 13).Value2
    P = $ws.Cells.Item(12, 16).Value2
}
$orig[13] = @{
    D = $ws.Cells.Item(13, 4).Value2
    J = $ws.Cells.Item(13, 10).Value2
    K = $ws.Cells.Item(13, 11).Value2
    L = $ws.Cells.Item(13, 12).Value2
    M = $ws.Cells.Item(13, 13).Value2
    P = $ws.Cells.Item(13, 16).Value2
}
$orig[14] = @{
    D = $ws.Cells.Item(14, 4).Value2
    J = $ws.Cells.Item(14, 10).Value2
    K = $ws.Cells.Item(14, 11).Value2
    L = $ws.Cells.Item(14, 12).Value2
    M = $ws.Cells.Item(14, 13).Value2
    P = $ws.Cells.Item(14, 16).Value2
}
$orig[15] = @{
    D = $ws.Cells.Item(15, 4).Value2
    J = $ws.Cells.Item(15, 10).Value2
    K = $ws.Cells.Item(15, 11).Value2
    L = $ws.Cells.Item(15, 12).Value2
    M = $ws.Cells.Item(15, 13).Value2
    P = $ws.Cells.Item(15, 16).Value2
}
$orig[16] = @{
    D = $ws.Cells.Item(16, 4).Value2
    J = $ws.Cells.Item(16, 10).Value2
    K = $ws.Cells.Item(16, 11).Value2
    L = $ws.Cells.Item(16, 12).Value2
    M = $ws.Cells.Item(16, 13).Value2
    P = $ws.Cells.Item(16, 16).Value2
}
$orig[17] = @{
    D = $ws.Cells.Item(17, 4).Value2
    J = $ws.Cells.Item(17, 10).Value2
    K = $ws.Cells.Item(17, 11).Value2
    L = $ws.Cells.Item(17, 12).Value2
    M = $ws.Cells.Item(17, 13).Value2
    P = $ws.Cells.Item(17, 16).Value2
}
$orig[18] = @{
    D = $ws.Cells.Item(18, 4).Value2
    J = $ws.Cells.Item(18, 10).Value2
    K = $ws.Cells.Item(18, 11).Value2
    L = $ws.Cells.Item(18, 12).Value2
    M = $ws.Cells.Item(18, 13).Value2
    P = $ws.Cells.Item(18, 16).Value2
}
$orig[20] = @{
    D = $ws.Cells.Item(20, 4).Value2
    J = $ws.Cells.Item(20, 10).Value2
    K = $ws.Cells.Item(20, 11).Value2
    L = $ws.Cells.Item(20, 12).Value2
    M = $ws.Cells.Item(20, 13).Value2
    P = $ws.Cells.Item(20, 16).Value2
}
$orig[21] = @{
    D = $ws.Cells.Item(21, 4).Value2
    J = $ws.Cells.Item(21, 10).Value2
    K = $ws.Cells.Item(21, 11).Value2
    L = $ws.Cells.Item(21, 12).Value2
    M = $ws.Cells.Item(21, 13).Value2
    P = $ws.Cells.Item(21, 16).Value2
}
$orig[22] = @{
    D = $ws.Cells.Item(22, 4).Value2
    J = $ws.Cells.Item(22, 10).Value2
    K = $ws.Cells.Item(22, 11).Value2
    L = $ws.Cells.Item(22, 12).Value2
    M = $ws.Cells.Item(22, 13).Value2
    P = $ws.Cells.Item(22, 16).Value2
}
$orig[23] = @{
    D = $ws.Cells.Item(23, 4).Value2
    J = $ws.Cells.Item(23, 10).Value2
    K = $ws.Cells.Item(23, 11).Value2
    L = $ws.Cells.Item(23, 12).Value2
    M = $ws.Cells.Item(23, 13).Value2
    P = $ws.Cells.Item(23, 16).Value2
}
$orig[24] = @{
    D = $ws.Cells.Item(24, 4).Value2
    J = $ws.Cells.Item(24, 10).Value2
    K = $ws.Cells.Item(24, 11).Value2
    L = $ws.Cells.Item(24, 12).Value2
    M = $ws.Cells.Item(24, 13).Value2
    P = $ws.Cells.Item(24, 16).Value2
}
$orig[25] = @{
    D = $ws.Cells.Item(25, 4).Value2
    J = $ws.Cells.Item(25, 10).Value2
    K = $ws.Cells.Item(25, 11).Value2
    L = $ws.Cells.Item(25, 12).Value2
    M = $ws.Cells.Item(25, 13).Value2
    P = $ws.Cells.Item(25, 16).Value2
}
$orig[26] = @{
    D = $ws.Cells.Item(26, 4).Value2
    J = $ws.Cells.Item(26, 10).Value2
    K = $ws.Cells.Item(26, 11).Value2
    L = $ws.Cells.Item(26, 12).Value2
    M = $ws.Cells.Item(26, 13).Value2
    P = $ws.Cells.Item(26, 16).Value2
}
$orig[27] = @{
    D = $ws.Cells.Item(27, 4).Value2
    J = $ws.Cells.Item(27, 10).Value2
    K = $ws.Cells.Item(27, 11).Value2
    L = $ws.Cells.Item(27, 12).Value2
    M = $ws.Cells.Item(27, 13).Value2
    P = $ws.Cells.Item(27, 16).Value2
}
$orig[28] = @{
    D = $ws.Cells.Item(28, 4).Value2
    J = $ws.Cells.Item(28, 10).Value2
    K = $ws.Cells.Item(28, 11).Value2
    L = $ws.Cells.Item(28, 12).Value2
    M = $ws.Cells.Item(28, 13).Value2
    P = $ws.Cells.Item(28, 16).Value2
}

# Apply permuted values: row N gets the original values of row mapping[N]
$ws.Cells.Item(2, 4).Value2 = $orig[15].D
$ws.Cells.Item(2, 10).Value2 = $orig[15].J
$ws.Cells.Item(2, 11).Value2 = $orig[15].K
$ws.Cells.Item(2, 12).Value2 = $orig[15].L
$ws.Cells.Item(2, 13).Value2 = $orig[15].M
$ws.Cells.Item(2, 16).Value2 = $orig[15].P

$ws.Cells.Item(3, 4).Value2 = $orig[24].D
$ws.Cells.Item(3, 10).Value2 = $orig[24].J
$ws.Cells.Item(3, 11).Value2 = $orig[24].K
$ws.Cells.Item(3, 12).Value2 = $orig[24].L
$ws.Cells.Item(3, 13).Value2 = $orig[24].M
$ws.Cells.Item(3, 16).Value2 = $orig[24].P

$ws.Cells.Item(4, 4).Value2 = $orig[13].D
$ws.Cells.Item(4, 10).Value2 = $orig[13].J
$ws.Cells.Item(4, 11).Value2 = $orig[13].K
$ws.Cells.Item(4, 12).Value2 = $orig[13].L
$ws.Cells.Item(4, 13).Value2 = $orig[13].M
$ws.Cells.Item(4, 16).Value2 = $orig[13].P

$ws.Cells.Item(5, 4).Value2 = $orig[22].D
$ws.Cells.Item(5, 10).Value2 = $orig[22].J
$ws.Cells.Item(5, 11).Value2 = $orig[22].K
$ws.Cells.Item(5, 12).Value2 = $orig[22].L
$ws.Cells.Item(5, 13).Value2 = $orig[22].M
$ws.Cells.Item(5, 16).Value2 = $orig[22].P

$ws.Cells.Item(6, 4).Value2 = $orig[26].D
$ws.Cells.Item(6, 10).Value2 = $orig[26].J
$ws.Cells.Item(6, 11).Value2 = $orig[26].K
$ws.Cells.Item(6, 12).Value2 = $orig[26].L
$ws.Cells.Item(6, 13).Value2 = $orig[26].M
$ws.Cells.Item(6, 16).Value2 = $orig[26].P

$ws.Cells.Item(7, 4).Value2 = $orig[17].D
$ws.Cells.Item(7, 10).Value2 = $orig[17].J
$ws.Cells.Item(7, 11).Value2 = $orig[17].K
$ws.Cells.Item(7, 12).Value2 = $orig[17].L
$ws.Cells.Item(7, 13).Value2 = $orig[17].M
$ws.Cells.Item(7, 16).Value2 = $orig[17].P

$ws.Cells.Item(8, 4).Value2 = $orig[18].D
$ws.Cells.Item(8, 10).Value2 = $orig[18].J
$ws.Cells.Item(8, 11).Value2 = $orig[18].K
$ws.Cells.Item(8, 12).Value2 = $orig[18].L
$ws.Cells.Item(8, 13).Value2 = $orig[18].M
$ws.Cells.Item(8, 16).Value2 = $orig[18].P

$ws.Cells.Item(9, 4).Value2 = $orig[4].D
$ws.Cells.Item(9, 10).Value2 = $orig[4].J
$ws.Cells.Item(9, 11).Value2 = $orig[4].K
$ws.Cells.Item(9, 12).Value2 = $orig[4].L
$ws.Cells.Item(9, 13).Value2 = $orig[4].M
$ws.Cells.Item(9, 16).Value2 = $orig[4].P

$ws.Cells.Item(10, 4).Value2 = $orig[14].D
$ws.Cells.Item(10, 10).Value2 = $orig[14].J
$ws.Cells.Item(10, 11).Value2 = $orig[14].K
$ws.Cells.Item(10, 12).Value2 = $orig[14].L
$ws.Cells.Item(10, 13).Value2 = $orig[14].M
$ws.Cells.Item(10, 16).Value2 = $orig[14].P

$ws.Cells.Item(11, 4).Value2 = $orig[16].D
$ws.Cells.Item(11, 10).Value2 = $orig[16].J
$ws.Cells.Item(11, 11).Value2 = $orig[16].K
$ws.Cells.Item(11, 12).Value2 = $orig[16].L
$ws.Cells.Item(11, 13).Value2 = $orig[16].M
$ws.Cells.Item(11, 16).Value2 = $orig[16].P

$ws.Cells.Item(12, 4).Value2 = $orig[6].D
$ws.Cells.Item(12, 10).Value2 = $orig[6].J
$ws.Cells.Item(12, 11).Value2 = $orig[6].K
$ws.Cells.Item(12, 12).Value2 = $orig[6].L
$ws.Cells.Item(12, 13).Value2 = $orig[6].M
$ws.Cells.Item(12, 16).Value2 = $orig[6].P

$ws.Cells.Item(13, 4).Value2 = $orig[25].D
$ws.Cells.Item(13, 10).Value2 = $orig[25].J
$ws.Cells.Item(13, 11).Value2 = $orig[25].K
$ws.Cells.Item(13, 12).Value2 = $orig[25].L
$ws.Cells.Item(13, 13).Value2 = $orig[25].M
$ws.Cells.Item(13, 16).Value2 = $orig[25].P

$ws.Cells.Item(14, 4).Value2 = $orig[7].D
$ws.Cells.Item(14, 10).Value2 = $orig[7].J
$ws.Cells.Item(14, 11).Value2 = $orig[7].K
$ws.Cells.Item(14, 12).Value2 = $orig[7].L
$ws.Cells.Item(14, 13).Value2 = $orig[7].M
$ws.Cells.Item(14, 16).Value2 = $orig[7].P

$ws.Cells.Item(15, 4).Value2 = $orig[3].D
$ws.Cells.Item(15, 10).Value2 = $orig[3].J
$ws.Cells.Item(15, 11).Value2 = $orig[3].K
$ws.Cells.Item(15, 12).Value2 = $orig[3].L
$ws.Cells.Item(15, 13).Value2 = $orig[3].M
$ws.Cells.Item(15, 16).Value2 = $orig[3].P

$ws.Cells.Item(16, 4).Value2 = $orig[11].D
$ws.Cells.Item(16, 10).Value2 = $orig[11].J
$ws.Cells.Item(16, 11).Value2 = $orig[11].K
$ws.Cells.Item(16, 12).Value2 = $orig[11].L
$ws.Cells.Item(16, 13).Value2 = $orig[11].M
$ws.Cells.Item(16, 16).Value2 = $orig[11].P

$ws.Cells.Item(17, 4).Value2 = $orig[9].D
$ws.Cells.Item(17, 10).Value2 = $orig[9].J
$ws.Cells.Item(17, 11).Value2 = $orig[9].K
$ws.Cells.Item(17, 12).Value2 = $orig[9].L
$ws.Cells.Item(17, 13).Value2 = $orig[9].M
$ws.Cells.Item(17, 16).Value2 = $orig[9].P

$ws.Cells.Item(18, 4).Value2 = $orig[23].D
$ws.Cells.Item(18, 10).Value2 = $orig[23].J
$ws.Cells.Item(18, 11).Value2 = $orig[23].K
$ws.Cells.Item(18, 12).Value2 = $orig[23].L
$ws.Cells.Item(18, 13).Value2 = $orig[23].M
$ws.Cells.Item(18, 16).Value2 = $orig[23].P

$ws.Cells.Item(20, 4).Value2 = $orig[10].D
$ws.Cells.Item(20, 10).Value2 = $orig[10].J
$ws.Cells.Item(20, 11).Value2 = $orig[10].K
$ws.Cells.Item(20, 12).Value2 = $orig[10].L
$ws.Cells.Item(20, 13).Value2 = $orig[10].M
$ws.Cells.Item(20, 16).Value2 = $orig[10].P

$ws.Cells.Item(21, 4).Value2 = $orig[2].D
$ws.Cells.Item(21, 10).Value2 = $orig[2].J
$ws.Cells.Item(21, 11).Value2 = $orig[2].K
$ws.Cells.Item(21, 12).Value2 = $orig[2].L
$ws.Cells.Item(21, 13).Value2 = $orig[2].M
$ws.Cells.Item(21, 16).Value2 = $orig[2].P

$ws.Cells.Item(22, 4).Value2 = $orig[28].D
$ws.Cells.Item(22, 10).Value2 = $orig[28].J
$ws.Cells.Item(22, 11).Value2 = $orig[28].K
$ws.Cells.Item(22, 12).Value2 = $orig[28].L
$ws.Cells.Item(22, 13).Value2 = $orig[28].M
$ws.Cells.Item(22, 16).Value2 = $orig[28].P

$ws.Cells.Item(23, 4).Value2 = $orig[5].D
$ws.Cells.Item(23, 10).Value2 = $orig[5].J
$ws.Cells.Item(23, 11).Value2 = $orig[5].K
$ws.Cells.Item(23, 12).Value2 = $orig[5].L
$ws.Cells.Item(23, 13).Value2 = $orig[5].M
$ws.Cells.Item(23, 16).Value2 = $orig[5].P

$ws.Cells.Item(24, 4).Value2 = $orig[20].D
$ws.Cells.Item(24, 10).Value2 = $orig[20].J
$ws.Cells.Item(24, 11).Value2 = $orig[20].K
$ws.Cells.Item(24, 12).Value2 = $orig[20].L
$ws.Cells.Item(24, 13).Value2 = $orig[20].M
$ws.Cells.Item(24, 16).Value2 = $orig[20].P

$ws.Cells.Item(25, 4).Value2 = $orig[12].D
$ws.Cells.Item(25, 10).Value2 = $orig[12].J
$ws.Cells.Item(25, 11).Value2 = $orig[12].K
$ws.Cells.Item(25, 12).Value2 = $orig[12].L
$ws.Cells.Item(25, 13).Value2 = $orig[12].M
$ws.Cells.Item(25, 16).Value2 = $orig[12].P

$ws.Cells.Item(26, 4).Value2 = $orig[27].D
$ws.Cells.Item(26, 10).Value2 = $orig[27].J
$ws.Cells.Item(26, 11).Value2 = $orig[27].K
$ws.Cells.Item(26, 12).Value2 = $orig[27].L
$ws.Cells.Item(26, 13).Value2 = $orig[27].M
$ws.Cells.Item(26, 16).Value2 = $orig[27].P

$ws.Cells.Item(27, 4).Value2 = $orig[8].D
$ws.Cells.Item(27, 10).Value2 = $orig[8].J
$ws.Cells.Item(27, 11).Value2 = $orig[8].K
$ws.Cells.Item(27, 12).Value2 = $orig[8].L
$ws.Cells.Item(27, 13).Value2 = $orig[8].M
$ws.Cells.Item(27, 16).Value2 = $orig[8].P

$ws.Cells.Item(28, 4).Value2 = $orig[21].D
$ws.Cells.Item(28, 10).Value2 = $orig[21].J
$ws.Cells.Item(28, 11).Value2 = $orig[21].K
$ws.Cells.Item(28, 12).Value2 = $orig[21].L
$ws.Cells.Item(28, 13).Value2 = $orig[21].M
$ws.Cells.Item(28, 16).Value2 = $orig[21].P
